# Auto-generated edit script: updates TPM-derived values in Efna1-Epha1.xlsx
# per commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 21.23829066666667
$ws.Range("H2").Value = 63.714872
$ws.Range("I2").Value = 0.9042366413687101
$ws.Range("J2").Value = 0.90423664136871
$ws.Range("M2").Value = 3.086228
$ws.Range("N2").Value = 9.258683999999999
$ws.Range("O2").Value = 0.1108770174431943
$ws.Range("P2").Value = 0.1108770174431943
$ws.Range("Q2").Value = 65.54620732760533
$ws.Range("R2").Value = 589.915865948448
$ws.Range("S2").Value = 0.1002590618578139
$ws.Range("T2").Value = 0.1002590618578139
# Row 3
$ws.Range("G3").Value = 21.23829066666667
$ws.Range("H3").Value = 63.714872
$ws.Range("I3").Value = 0.9042366413687101
$ws.Range("J3").Value = 0.90423664136871
$ws.Range("M3").Value = 5.619945333333334
$ws.Range("O3").Value = 0.2019043235800461
$ws.Range("P3").Value = 0.2019043235800461
$ws.Range("Q3").Value = 119.3580325201102
$ws.Range("R3").Value = 1074.222292680992
$ws.Range("S3").Value = 0.1825692874318422
$ws.Range("T3").Value = 0.1825692874318422
# Row 4
$ws.Range("G4").Value = 21.23829066666667
$ws.Range("H4").Value = 63.714872
$ws.Range("I4").Value = 0.9042366413687101
$ws.Range("J4").Value = 0.90423664136871
$ws.Range("M4").Value = 3.355234666666667
$ws.Range("N4").Value = 10.065704
$ws.Range("O4").Value = 0.1205414547019891
$ws.Range("P4").Value = 0.1205414547019891
$ws.Range("Q4").Value = 71.25944910554311
$ws.Range("R4").Value = 641.335041949888
$ws.Range("S4").Value = 0.1089980001454251
$ws.Range("T4").Value = 0.1089980001454251
# Row 5
$ws.Range("G5").Value = 21.23829066666667
$ws.Range("H5").Value = 63.714872
$ws.Range("I5").Value = 0.9042366413687101
$ws.Range("J5").Value = 0.90423664136871
$ws.Range("M5").Value = 1.012758666666667
$ws.Range("N5").Value = 3.038276
$ws.Range("O5").Value = 0.03638475846559173
$ws.Range("P5").Value = 0.03638475846559173
$ws.Range("Q5").Value = 21.50926293785244
$ws.Range("R5").Value = 193.583366440672
$ws.Range("S5").Value = 0.0329004317919384
$ws.Range("T5").Value = 0.0329004317919384
# Row 6
$ws.Range("G6").Value = 21.23829066666667
$ws.Range("H6").Value = 63.714872
$ws.Range("I6").Value = 0.9042366413687101
$ws.Range("J6").Value = 0.90423664136871
$ws.Range("M6").Value = 12.86621566666667
$ws.Range("N6").Value = 38.598647
$ws.Range("O6").Value = 0.4622366263610143
$ws.Range("P6").Value = 0.4622366263610142
$ws.Range("Q6").Value = 273.2564281086871
$ws.Range("R6").Value = 2459.307852978184
$ws.Range("S6").Value = 0.4179712945382869
$ws.Range("T6").Value = 0.4179712945382868
# Row 7
$ws.Range("G7").Value = 21.23829066666667
$ws.Range("H7").Value = 63.714872
$ws.Range("I7").Value = 0.9042366413687101
$ws.Range("J7").Value = 0.90423664136871
$ws.Range("M7").Value = 1.894313
$ws.Range("N7").Value = 5.682938999999999
$ws.Range("O7").Value = 0.06805581944816448
$ws.Range("P7").Value = 0.06805581944816448
$ws.Range("Q7").Value = 40.23197010764533
$ws.Range("R7").Value = 362.0877309688079
$ws.Range("S7").Value = 0.06153856560340359
$ws.Range("T7").Value = 0.06153856560340358
# Row 8
$ws.Range("I8").Value = 0.05937834432696559
$ws.Range("J8").Value = 0.05937834432696559
$ws.Range("M8").Value = 3.086228
$ws.Range("N8").Value = 9.258683999999999
$ws.Range("O8").Value = 0.1108770174431943
$ws.Range("P8").Value = 0.1108770174431943
$ws.Range("Q8").Value = 4.304210966428
$ws.Range("R8").Value = 38.737898697852
$ws.Range("S8").Value = 0.006583693719688964
$ws.Range("T8").Value = 0.006583693719688963
# Row 9
$ws.Range("I9").Value = 0.05937834432696559
$ws.Range("J9").Value = 0.05937834432696559
$ws.Range("M9").Value = 5.619945333333334
$ws.Range("O9").Value = 0.2019043235800461
$ws.Range("P9").Value = 0.2019043235800461
$ws.Range("Q9").Value = 7.837862379078669
$ws.Range("R9").Value = 70.54076141170802
$ws.Range("S9").Value = 0.01198874444663906
$ws.Range("T9").Value = 0.01198874444663906
# Row 10
$ws.Range("I10").Value = 0.05937834432696559
$ws.Range("J10").Value = 0.05937834432696559
$ws.Range("M10").Value = 3.355234666666667
$ws.Range("N10").Value = 10.065704
$ws.Range("O10").Value = 0.1205414547019891
$ws.Range("P10").Value = 0.1205414547019891
$ws.Range("Q10").Value = 4.679381383101334
$ws.Range("R10").Value = 42.11443244791201
$ws.Range("S10").Value = 0.007157552002968033
$ws.Range("T10").Value = 0.007157552002968033
# Row 11
$ws.Range("I11").Value = 0.05937834432696559
$ws.Range("J11").Value = 0.05937834432696559
$ws.Range("M11").Value = 1.012758666666667
$ws.Range("N11").Value = 3.038276
$ws.Range("O11").Value = 0.03638475846559173
$ws.Range("P11").Value = 0.03638475846559173
$ws.Range("Q11").Value = 1.412444887225333
$ws.Range("R11").Value = 12.712003985028
$ws.Range("S11").Value = 0.002160466716423382
$ws.Range("T11").Value = 0.002160466716423382
# Row 12
$ws.Range("I12").Value = 0.05937834432696559
$ws.Range("J12").Value = 0.05937834432696559
$ws.Range("M12").Value = 12.86621566666667
$ws.Range("N12").Value = 38.598647
$ws.Range("O12").Value = 0.4622366263610143
$ws.Range("P12").Value = 0.4622366263610142
$ws.Range("Q12").Value = 17.94388054573234
$ws.Range("R12").Value = 161.494924911591
$ws.Range("S12").Value = 0.02744684556059925
$ws.Range("T12").Value = 0.02744684556059924
# Row 13
$ws.Range("I13").Value = 0.05937834432696559
$ws.Range("J13").Value = 0.05937834432696559
$ws.Range("M13").Value = 1.894313
$ws.Range("N13").Value = 5.682938999999999
$ws.Range("O13").Value = 0.06805581944816448
$ws.Range("P13").Value = 0.06805581944816448
$ws.Range("Q13").Value = 2.641905519763
$ws.Range("R13").Value = 23.777149677867
$ws.Range("S13").Value = 0.004041041880646912
$ws.Range("T13").Value = 0.004041041880646912
# Row 14
$ws.Range("G14").Value = 0.6246503333333333
$ws.Range("H14").Value = 1.873951
$ws.Range("I14").Value = 0.02659497076804196
$ws.Range("J14").Value = 0.02659497076804196
$ws.Range("M14").Value = 3.086228
$ws.Range("N14").Value = 9.258683999999999
$ws.Range("O14").Value = 0.1108770174431943
$ws.Range("P14").Value = 0.1108770174431943
$ws.Range("Q14").Value = 1.927813348942667
$ws.Range("R14").Value = 17.350320140484
$ws.Range("S14").Value = 0.002948771037749432
$ws.Range("T14").Value = 0.002948771037749432
# Row 15
$ws.Range("G15").Value = 0.6246503333333333
$ws.Range("H15").Value = 1.873951
$ws.Range("I15").Value = 0.02659497076804196
$ws.Range("J15").Value = 0.02659497076804196
$ws.Range("M15").Value = 5.619945333333334
$ws.Range("O15").Value = 0.2019043235800461
$ws.Range("P15").Value = 0.2019043235800461
$ws.Range("Q15").Value = 3.510500725781778
$ws.Range("R15").Value = 31.594506532036
$ws.Range("S15").Value = 0.005369639583552613
$ws.Range("T15").Value = 0.005369639583552613
# Row 16
$ws.Range("G16").Value = 0.6246503333333333
$ws.Range("H16").Value = 1.873951
$ws.Range("I16").Value = 0.02659497076804196
$ws.Range("J16").Value = 0.02659497076804196
$ws.Range("M16").Value = 3.355234666666667
$ws.Range("N16").Value = 10.065704
$ws.Range("O16").Value = 0.1205414547019891
$ws.Range("P16").Value = 0.1205414547019891
$ws.Range("Q16").Value = 2.095848452944889
$ws.Range("R16").Value = 18.862636076504
$ws.Range("S16").Value = 0.003205796464136653
$ws.Range("T16").Value = 0.003205796464136653
# Row 17
$ws.Range("G17").Value = 0.6246503333333333
$ws.Range("H17").Value = 1.873951
$ws.Range("I17").Value = 0.02659497076804196
$ws.Range("J17").Value = 0.02659497076804196
$ws.Range("M17").Value = 1.012758666666667
$ws.Range("N17").Value = 3.038276
$ws.Range("O17").Value = 0.03638475846559173
$ws.Range("P17").Value = 0.03638475846559173
$ws.Range("Q17").Value = 0.6326200387195555
$ws.Range("R17").Value = 5.693580348475999
$ws.Range("S17").Value = 0.0009676515877946793
$ws.Range("T17").Value = 0.0009676515877946793
# Row 18
$ws.Range("G18").Value = 0.6246503333333333
$ws.Range("H18").Value = 1.873951
$ws.Range("I18").Value = 0.02659497076804196
$ws.Range("J18").Value = 0.02659497076804196
$ws.Range("M18").Value = 12.86621566666667
$ws.Range("N18").Value = 38.598647
$ws.Range("O18").Value = 0.4622366263610143
$ws.Range("P18").Value = 0.4622366263610142
$ws.Range("Q18").Value = 8.036885904921888
$ws.Range("R18").Value = 72.331973144297
$ws.Range("S18").Value = 0.01229316956598951
$ws.Range("T18").Value = 0.01229316956598951
# Row 19
$ws.Range("G19").Value = 0.6246503333333333
$ws.Range("H19").Value = 1.873951
$ws.Range("I19").Value = 0.02659497076804196
$ws.Range("J19").Value = 0.02659497076804196
$ws.Range("M19").Value = 1.894313
$ws.Range("N19").Value = 5.682938999999999
$ws.Range("O19").Value = 0.06805581944816448
$ws.Range("P19").Value = 0.06805581944816448
$ws.Range("Q19").Value = 1.183283246887666
$ws.Range("R19").Value = 10.649549221989
$ws.Range("S19").Value = 0.001809942528819076
$ws.Range("T19").Value = 0.001809942528819076
# Row 20
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 0.6666666666666666
$ws.Range("G20").Value = 0.229944
$ws.Range("H20").Value = 0.689832
$ws.Range("I20").Value = 0.009790043536282392
$ws.Range("J20").Value = 0.009790043536282392
$ws.Range("M20").Value = 3.086228
$ws.Range("N20").Value = 9.258683999999999
$ws.Range("O20").Value = 0.1108770174431943
$ws.Range("P20").Value = 0.1108770174431943
$ws.Range("Q20").Value = 0.7096596112319999
$ws.Range("R20").Value = 6.386936501087999
$ws.Range("S20").Value = 0.001085490827942015
$ws.Range("T20").Value = 0.001085490827942015
# Row 21
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 0.6666666666666666
$ws.Range("G21").Value = 0.229944
$ws.Range("H21").Value = 0.689832
$ws.Range("I21").Value = 0.009790043536282392
$ws.Range("J21").Value = 0.009790043536282392
$ws.Range("M21").Value = 5.619945333333334
$ws.Range("O21").Value = 0.2019043235800461
$ws.Range("P21").Value = 0.2019043235800461
$ws.Range("Q21").Value = 1.292272709728
$ws.Range("R21").Value = 11.630454387552
$ws.Range("S21").Value = 0.001976652118012299
$ws.Range("T21").Value = 0.001976652118012299
# Row 22
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 0.6666666666666666
$ws.Range("G22").Value = 0.229944
$ws.Range("H22").Value = 0.689832
$ws.Range("I22").Value = 0.009790043536282392
$ws.Range("J22").Value = 0.009790043536282392
$ws.Range("M22").Value = 3.355234666666667
$ws.Range("N22").Value = 10.065704
$ws.Range("O22").Value = 0.1205414547019891
$ws.Range("P22").Value = 0.1205414547019891
$ws.Range("Q22").Value = 0.771516080192
$ws.Range("R22").Value = 6.943644721728
$ws.Range("S22").Value = 0.001180106089459285
$ws.Range("T22").Value = 0.001180106089459285
# Row 23
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 0.6666666666666666
$ws.Range("G23").Value = 0.229944
$ws.Range("H23").Value = 0.689832
$ws.Range("I23").Value = 0.009790043536282392
$ws.Range("J23").Value = 0.009790043536282392
$ws.Range("M23").Value = 1.012758666666667
$ws.Range("N23").Value = 3.038276
$ws.Range("O23").Value = 0.03638475846559173
$ws.Range("P23").Value = 0.03638475846559173
$ws.Range("Q23").Value = 0.232877778848
$ws.Range("R23").Value = 2.095900009632
$ws.Range("S23").Value = 0.0003562083694352623
$ws.Range("T23").Value = 0.0003562083694352623
# Row 24
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 0.6666666666666666
$ws.Range("G24").Value = 0.229944
$ws.Range("H24").Value = 0.689832
$ws.Range("I24").Value = 0.009790043536282392
$ws.Range("J24").Value = 0.009790043536282392
$ws.Range("M24").Value = 12.86621566666667
$ws.Range("N24").Value = 38.598647
$ws.Range("O24").Value = 0.4622366263610143
$ws.Range("P24").Value = 0.4622366263610142
$ws.Range("Q24").Value = 2.958509095256
$ws.Range("R24").Value = 26.626581857304
$ws.Range("S24").Value = 0.004525316696138626
$ws.Range("T24").Value = 0.004525316696138626
# Row 25
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 0.6666666666666666
$ws.Range("G25").Value = 0.229944
$ws.Range("H25").Value = 0.689832
$ws.Range("I25").Value = 0.009790043536282392
$ws.Range("J25").Value = 0.009790043536282392
$ws.Range("M25").Value = 1.894313
$ws.Range("N25").Value = 5.682938999999999
$ws.Range("O25").Value = 0.06805581944816448
$ws.Range("P25").Value = 0.06805581944816448
$ws.Range("Q25").Value = 0.435585908472
$ws.Range("R25").Value = 3.920273176248
$ws.Range("S25").Value = 0.0006662694352949042
$ws.Range("T25").Value = 0.0006662694352949042
